$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "245.23"
$ws.Range("D3").Value = "23.07"
$ws.Range("D4").Value = "5.409"
$ws.Range("D5").Value = "0.06003"
$ws.Range("D8").Value = "0.9277"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "0.0005943"
$ws.Range("E9").Value = "8OneONE"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1419"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07430"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03385"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03044"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09357"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "3.936"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001589"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04803"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("D18").Value = "0.005325"
$ws.Range("D19").Value = "0.004156"
$ws.Range("D20").Value = "0.0009823"
$ws.Range("D22").Value = "3.664"
$ws.Range("D23").Value = "6.425"
$ws.Range("D24").Value = "2.187"
$ws.Range("D40").Value = "0.03968"
$ws.Range("D41").Value = "0.006394"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("D43").Value = "0.002902"
$ws.Range("D44").Value = "0.006806"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("D45").Value = "0.00005202"
$ws.Range("D47").Value = "0.0005803"
$ws.Range("D48").Value = "0.9806"
